# Updates computed market-price / profit columns (H-N) across the ALC, ARM, BSM, CRP,
# CUL, GSM, LTW and WVR leve-profit sheets, refreshing cached marketboard prices.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 333590
$ws.Range("I9").Value = 1000000
$ws.Range("J9").Value = 385
$ws.Range("K9").Value = 1000000
$ws.Range("L9").Value = 385
$ws.Range("M9").Value = -999831
$ws.Range("N9").Value = -723
# Row 43
$ws.Range("H43").Value = 966.6667
$ws.Range("I43").Value = 900
$ws.Range("J43").Value = 1000
$ws.Range("K43").Value = 900
$ws.Range("L43").Value = 1000
$ws.Range("M43").Value = -831
$ws.Range("N43").Value = -1138
# Row 86
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
# Row 89
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
# Row 92
$ws.Range("H92").Value = 1191.6364
$ws.Range("I92").Value = 1156.7222
$ws.Range("K92").Value = 1156.7222
$ws.Range("M92").Value = 91.27780000000007
# Row 107
$ws.Range("H107").Value = 1283.5
$ws.Range("I107").Value = 1241.3529
$ws.Range("K107").Value = 1241.3529
$ws.Range("M107").Value = 678.6470999999999
# Row 129
$ws.Range("H129").Value = 1994.45
$ws.Range("I129").Value = 1497.5
$ws.Range("J129").Value = 2118.6875
$ws.Range("K129").Value = 4492.5
$ws.Range("L129").Value = 6356.0625
$ws.Range("M129").Value = 507.5
$ws.Range("N129").Value = -16356.0625
# Row 132
$ws.Range("H132").Value = 1878.8
$ws.Range("I132").Value = 1878.8
$ws.Range("K132").Value = 5636.4
$ws.Range("M132").Value = -3106.4

$ws = $wb.Worksheets.Item("ARM")
# Row 46
$ws.Range("H46").Value = 14785
$ws.Range("I46").Value = 10000
$ws.Range("K46").Value = 10000
$ws.Range("M46").Value = -9681
# Row 74
$ws.Range("H74").Value = 18146.285
$ws.Range("I74").Value = 13400.6
$ws.Range("K74").Value = 13400.6
$ws.Range("M74").Value = -12526.6
# Row 77
$ws.Range("H77").Value = 18146.285
$ws.Range("I77").Value = 13400.6
$ws.Range("K77").Value = 67003
$ws.Range("M77").Value = -62635

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 8725.556
$ws.Range("J20").Value = 7605.3335
$ws.Range("L20").Value = 7605.3335
$ws.Range("N20").Value = -8099.3335
# Row 81
$ws.Range("H81").Value = 40000
$ws.Range("J81").Value = 40000
$ws.Range("L81").Value = 40000
$ws.Range("N81").Value = -42122
# Row 84
$ws.Range("H84").Value = 40000
$ws.Range("J84").Value = 40000
$ws.Range("L84").Value = 120000
$ws.Range("N84").Value = -130608
# Row 86
$ws.Range("H86").Value = 5539.778
$ws.Range("I86").Value = 3107.25
$ws.Range("K86").Value = 3107.25
$ws.Range("M86").Value = -1984.25
# Row 89
$ws.Range("H89").Value = 5539.778
$ws.Range("I89").Value = 3107.25
$ws.Range("K89").Value = 15536.25
$ws.Range("M89").Value = -9920.25
# Row 94
$ws.Range("H94").Value = 4575.636
$ws.Range("I94").Value = 4583.2
$ws.Range("J94").Value = 4500
$ws.Range("K94").Value = 4583.2
$ws.Range("L94").Value = 4500
$ws.Range("M94").Value = -4132.2
$ws.Range("N94").Value = -5402
# Row 134
$ws.Range("H134").Value = 1101.2
$ws.Range("I134").Value = 1101.2
$ws.Range("K134").Value = 3303.6
$ws.Range("M134").Value = -768.6000000000004

$ws = $wb.Worksheets.Item("CRP")
# Row 132
$ws.Range("H132").Value = 3897
$ws.Range("I132").Value = 3344
$ws.Range("K132").Value = 10032
$ws.Range("M132").Value = -7502

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 16515
$ws.Range("I3").Value = 16515
$ws.Range("K3").Value = 49545
$ws.Range("M3").Value = -49433
# Row 101
$ws.Range("H101").Value = 7857.143
$ws.Range("J101").Value = 7857.143
$ws.Range("L101").Value = 23571.429
$ws.Range("N101").Value = -28439.429
# Row 113
$ws.Range("H113").Value = 2166.3333
$ws.Range("J113").Value = 2166.3333
$ws.Range("L113").Value = 6498.999899999999
$ws.Range("N113").Value = -10838.9999
# Row 122
$ws.Range("H122").Value = 1784.4
$ws.Range("J122").Value = 2327.8333
$ws.Range("L122").Value = 20950.4997
$ws.Range("N122").Value = -25850.4997
# Row 132
$ws.Range("H132").Value = 2475.8823
$ws.Range("J132").Value = 3359.1
$ws.Range("L132").Value = 30231.9
$ws.Range("N132").Value = -35291.89999999999
# Row 136
$ws.Range("H136").Value = 2348.7058
$ws.Range("I136").Value = 2009.3334
$ws.Range("J136").Value = 2421.4285
$ws.Range("K136").Value = 6028.0002
$ws.Range("L136").Value = 7264.2855
$ws.Range("M136").Value = -928.0002000000004
$ws.Range("N136").Value = -17464.2855

$ws = $wb.Worksheets.Item("GSM")
# Row 86
$ws.Range("H86").Value = 40000
$ws.Range("J86").Value = 40000
$ws.Range("L86").Value = 40000
$ws.Range("N86").Value = -42372
# Row 89
$ws.Range("H89").Value = 40000
$ws.Range("J89").Value = 40000
$ws.Range("L89").Value = 120000
$ws.Range("N89").Value = -131856
# Row 113
$ws.Range("H113").Value = 1825.6666
$ws.Range("I113").Value = 1453.5
$ws.Range("K113").Value = 1453.5
$ws.Range("M113").Value = 716.5
# Row 132
$ws.Range("H132").Value = 3180.2354
$ws.Range("I132").Value = 1962.7
$ws.Range("K132").Value = 5888.1
$ws.Range("M132").Value = -3358.1

$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 2148.5454
$ws.Range("I82").Value = 2028.1666
$ws.Range("K82").Value = 2028.1666
$ws.Range("M82").Value = -1667.1666
# Row 85
$ws.Range("H85").Value = 2148.5454
$ws.Range("I85").Value = 2028.1666
$ws.Range("K85").Value = 2028.1666
$ws.Range("M85").Value = -780.1666
# Row 136
$ws.Range("H136").Value = 4732.8887
$ws.Range("I136").Value = 1765.3334
$ws.Range("K136").Value = 5296.0002
$ws.Range("M136").Value = -2746.0002

$ws = $wb.Worksheets.Item("WVR")
# Row 68
$ws.Range("H68").Value = 40000
$ws.Range("J68").Value = 40000
$ws.Range("L68").Value = 40000
$ws.Range("N68").Value = -41622
# Row 71
$ws.Range("H71").Value = 40000
$ws.Range("J71").Value = 40000
$ws.Range("L71").Value = 120000
$ws.Range("N71").Value = -128112
# Row 81
$ws.Range("H81").Value = 4840.3335
$ws.Range("I81").Value = 2760.5
$ws.Range("K81").Value = 5521
$ws.Range("M81").Value = -4460
# Row 84
$ws.Range("H84").Value = 4840.3335
$ws.Range("I84").Value = 2760.5
$ws.Range("K84").Value = 27605
$ws.Range("M84").Value = -22301
# Row 113
$ws.Range("H113").Value = 1117.7
$ws.Range("I113").Value = 1539.6
$ws.Range("K113").Value = 4618.799999999999
$ws.Range("M113").Value = -2448.799999999999
# Row 132
$ws.Range("H132").Value = 3529.353
$ws.Range("I132").Value = 1691.75
$ws.Range("K132").Value = 5075.25
$ws.Range("M132").Value = -2545.25
